$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correção das notas do fórum: zera as colunas B:J (dias de acesso,
# total_views e nota_view) para todas as linhas de alunos (2 a 50).
$ws.Range("B2:J50").Value = 0
